$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the newly-reviewed "Kész?" (Done?) column for the GET-parameter,
# complex-calculation and exception-handling rows.
# Order matters for how new shared strings are interned: "IGEN(B)" must be
# created before "IGEN(B) (almost)" to match the target workbook.
$ws.Range("D26").Value = "IGEN(B)"
$ws.Range("D23").Value = "IGEN(B) (almost)"
$ws.Range("D32").Value = "IGEN (B)"

# Scroll the sheet so row 7 is at the top and select D23, matching where the
# author was working when they saved.
$ws.Range("D23").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
